# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Pandaemonium_Profits leve-profit workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4124.5835
$ws.Range("I64").Value = 3649.2856
$ws.Range("K64").Value = 3649.2856
$ws.Range("M64").Value = -3401.2856

$ws.Range("H67").Value = 4124.5835
$ws.Range("I67").Value = 3649.2856
$ws.Range("K67").Value = 3649.2856
$ws.Range("M67").Value = -2791.2856

$ws.Range("H70").Value = 1625.421
$ws.Range("I70").Value = 1278.1
$ws.Range("J70").Value = 2011.3334
$ws.Range("K70").Value = 3834.3
$ws.Range("L70").Value = 6034.0002
$ws.Range("M70").Value = -3564.3
$ws.Range("N70").Value = -6574.0002

$ws.Range("H73").Value = 1625.421
$ws.Range("I73").Value = 1278.1
$ws.Range("J73").Value = 2011.3334
$ws.Range("K73").Value = 3834.3
$ws.Range("L73").Value = 6034.0002
$ws.Range("M73").Value = -2898.3
$ws.Range("N73").Value = -7906.0002

$ws.Range("H76").Value = 3908.6956
$ws.Range("I76").Value = 3730.9375
$ws.Range("K76").Value = 3730.9375
$ws.Range("M76").Value = -3415.9375

$ws.Range("H79").Value = 3908.6956
$ws.Range("I79").Value = 3730.9375
$ws.Range("K79").Value = 3730.9375
$ws.Range("M79").Value = -2638.9375

$ws.Range("H116").Value = 2577
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 3207.8
$ws.Range("K116").Value = 1000
$ws.Range("L116").Value = 3207.8
$ws.Range("M116").Value = 2442
$ws.Range("N116").Value = -10091.8

$ws.Range("H137").Value = 3668.8647
$ws.Range("I137").Value = 1838.409
$ws.Range("J137").Value = 6353.533
$ws.Range("K137").Value = 5515.227000000001
$ws.Range("L137").Value = 19060.599
$ws.Range("M137").Value = -2965.227000000001
$ws.Range("N137").Value = -24160.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9210.540999999999
$ws.Range("I61").Value = 4543.32
$ws.Range("J61").Value = 18933.916
$ws.Range("K61").Value = 4543.32
$ws.Range("L61").Value = 18933.916
$ws.Range("M61").Value = -4331.32
$ws.Range("N61").Value = -19357.916

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 5442.5366
$ws.Range("I132").Value = 2091.0588
$ws.Range("J132").Value = 7816.5
$ws.Range("K132").Value = 6273.176399999999
$ws.Range("L132").Value = 23449.5
$ws.Range("M132").Value = -3743.176399999999
$ws.Range("N132").Value = -28509.5

$ws.Range("H136").Value = 9210.540999999999
$ws.Range("I136").Value = 4543.32
$ws.Range("J136").Value = 18933.916
$ws.Range("K136").Value = 13629.96
$ws.Range("L136").Value = 56801.74800000001
$ws.Range("M136").Value = -11079.96
$ws.Range("N136").Value = -61901.74800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5531.276
$ws.Range("I105").Value = 5199.909
$ws.Range("K105").Value = 5199.909
$ws.Range("M105").Value = -3452.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2333853.5
$ws.Range("I58").Value = 3498675.2
$ws.Range("J58").Value = 4210.231
$ws.Range("K58").Value = 3498675.2
$ws.Range("L58").Value = 4210.231
$ws.Range("M58").Value = -3498472.2
$ws.Range("N58").Value = -4616.231

$ws.Range("H62").Value = 3101
$ws.Range("I62").Value = 3001.25
$ws.Range("K62").Value = 3001.25
$ws.Range("M62").Value = -2377.25

$ws.Range("H65").Value = 3101
$ws.Range("I65").Value = 3001.25
$ws.Range("K65").Value = 15006.25
$ws.Range("M65").Value = -11886.25

$ws.Range("H107").Value = 620.96
$ws.Range("I107").Value = 517.5294
$ws.Range("J107").Value = 840.75
$ws.Range("K107").Value = 517.5294
$ws.Range("L107").Value = 840.75
$ws.Range("M107").Value = 1402.4706
$ws.Range("N107").Value = -4680.75

$ws.Range("H136").Value = 2333853.5
$ws.Range("I136").Value = 3498675.2
$ws.Range("J136").Value = 4210.231
$ws.Range("K136").Value = 10496025.6
$ws.Range("L136").Value = 12630.693
$ws.Range("M136").Value = -10493475.6
$ws.Range("N136").Value = -17730.693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I130").Value = 940
$ws.Range("K130").Value = 2820
$ws.Range("M130").Value = 2200

$ws.Range("H131").Value = 18026.389
$ws.Range("I131").Value = 386.2245
$ws.Range("K131").Value = 1158.6735
$ws.Range("M131").Value = 3881.3265

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 43000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 43000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 43000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -44372

$ws.Range("H65").Value = 43000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 43000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 129000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -135864

$ws.Range("H80").Value = 10110.714
$ws.Range("I80").Value = 18326.666
$ws.Range("J80").Value = 3948.75
$ws.Range("K80").Value = 18326.666
$ws.Range("L80").Value = 3948.75
$ws.Range("M80").Value = -17328.666
$ws.Range("N80").Value = -5944.75

$ws.Range("H83").Value = 10110.714
$ws.Range("I83").Value = 18326.666
$ws.Range("J83").Value = 3948.75
$ws.Range("K83").Value = 91633.33
$ws.Range("L83").Value = 19743.75
$ws.Range("M83").Value = -86641.33
$ws.Range("N83").Value = -29727.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 35499.5
$ws.Range("J63").Value = 35499.5
$ws.Range("L63").Value = 35499.5
$ws.Range("N63").Value = -36997.5

$ws.Range("H66").Value = 35499.5
$ws.Range("J66").Value = 35499.5
$ws.Range("L66").Value = 106498.5
$ws.Range("N66").Value = -113986.5

$ws.Range("H132").Value = 3728.9443
$ws.Range("I132").Value = 3383.3635
$ws.Range("K132").Value = 10150.0905
$ws.Range("M132").Value = -7620.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11666.667
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 11666.667
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 11666.667
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -12446.667

$ws.Range("H45").Value = 4955.8
$ws.Range("I45").Value = 3593
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 3593
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = -3102
$ws.Range("N45").Value = -7982

$ws.Range("H108").Value = 48700
$ws.Range("J108").Value = 48700
$ws.Range("L108").Value = 48700
$ws.Range("N108").Value = -56380

$ws.Range("H126").Value = 1693.1428
$ws.Range("I126").Value = 1572
$ws.Range("J126").Value = 1814.2858
$ws.Range("K126").Value = 4716
$ws.Range("L126").Value = 5442.857400000001
$ws.Range("M126").Value = -2246
$ws.Range("N126").Value = -10382.8574

$ws.Range("H136").Value = 4216.4243
$ws.Range("I136").Value = 1734.7106
$ws.Range("J136").Value = 7584.4644
$ws.Range("K136").Value = 5204.1318
$ws.Range("L136").Value = 22753.3932
$ws.Range("M136").Value = -2654.1318
$ws.Range("N136").Value = -27853.3932
